# Update NATMI ligand-receptor pair statistics (Lgals3-Lag3) per Dr Hou advice.
# Columns E (Ligand-expressing cells) and K (Receptor-expressing cells) change from 1 to 3,
# and all dependent expression / specificity metrics (G-J, M-T) are recomputed accordingly.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$numRows = 16
$numCols = 16
$values = New-Object "object[,]" $numRows,$numCols

# Row 2: E,F,G,H,I,J,K,L,M,N,O,P,Q,R,S,T
$values[0,0] = 3
$values[0,1] = 1
$values[0,2] = 3.181738666666666
$values[0,3] = 9.545216
$values[0,4] = 0.01184997401866639
$values[0,5] = 0.01184997401866639
$values[0,6] = 3
$values[0,7] = 1
$values[0,8] = 10.38865533333333
$values[0,9] = 31.165966
$values[0,10] = 0.261404167660619
$values[0,11] = 0.261404167660619
$values[0,12] = 33.05398636873956
$values[0,13] = 297.485877318656
$values[0,14] = 0.003097632595149447
$values[0,15] = 0.003097632595149448

# Row 3: E,F,G,H,I,J,K,L,M,N,O,P,Q,R,S,T
$values[1,0] = 3
$values[1,1] = 1
$values[1,2] = 3.181738666666666
$values[1,3] = 9.545216
$values[1,4] = 0.01184997401866639
$values[1,5] = 0.01184997401866639
$values[1,6] = 3
$values[1,7] = 1
$values[1,8] = 10.56216766666667
$values[1,9] = 31.686503
$values[1,10] = 0.2657701655321932
$values[1,11] = 0.2657701655321932
$values[1,12] = 33.60605726884977
$values[1,13] = 302.454515419648
$values[1,14] = 0.003149369556493154
$values[1,15] = 0.003149369556493155

# Row 4: E,F,G,H,I,J,K,L,M,N,O,P,Q,R,S,T
$values[2,0] = 3
$values[2,1] = 1
$values[2,2] = 3.181738666666666
$values[2,3] = 9.545216
$values[2,4] = 0.01184997401866639
$values[2,5] = 0.01184997401866639
$values[2,6] = 3
$values[2,7] = 1
$values[2,8] = 14.90560833333333
$values[2,9] = 44.716825
$values[2,10] = 0.3750618357072762
$values[2,11] = 0.3750618357072762
$values[2,12] = 47.42575038435555
$values[2,13] = 426.8317534592
$values[2,14] = 0.004444473008524545
$values[2,15] = 0.004444473008524546

# Row 5: E,F,G,H,I,J,K,L,M,N,O,P,Q,R,S,T
$values[3,0] = 3
$values[3,1] = 1
$values[3,2] = 3.181738666666666
$values[3,3] = 9.545216
$values[3,4] = 0.01184997401866639
$values[3,5] = 0.01184997401866639
$values[3,6] = 3
$values[3,7] = 1
$values[3,8] = 3.885304333333333
$values[3,9] = 11.655913
$values[3,10] = 0.09776383109991163
$values[3,11] = 0.09776383109991162
$values[3,12] = 12.36202302913422
$values[3,13] = 111.258207262208
$values[3,14] = 0.001158498858499242
$values[3,15] = 0.001158498858499242

# Row 6: E,F,G,H,I,J,K,L,M,N,O,P,Q,R,S,T
$values[4,0] = 3
$values[4,1] = 1
$values[4,2] = 5.848171333333333
$values[4,3] = 17.544514
$values[4,4] = 0.02178075750932496
$values[4,5] = 0.02178075750932496
$values[4,6] = 3
$values[4,7] = 1
$values[4,8] = 10.38865533333333
$values[4,9] = 31.165966
$values[4,10] = 0.261404167660619
$values[4,11] = 0.261404167660619
$values[4,12] = 60.75463631228045
$values[4,13] = 546.791726810524
$values[4,14] = 0.005693580787742867
$values[4,15] = 0.005693580787742867

# Row 7: E,F,G,H,I,J,K,L,M,N,O,P,Q,R,S,T
$values[5,0] = 3
$values[5,1] = 1
$values[5,2] = 5.848171333333333
$values[5,3] = 17.544514
$values[5,4] = 0.02178075750932496
$values[5,5] = 0.02178075750932496
$values[5,6] = 3
$values[5,7] = 1
$values[5,8] = 10.56216766666667
$values[5,9] = 31.686503
$values[5,10] = 0.2657701655321932
$values[5,11] = 0.2657701655321932
$values[5,12] = 61.76936616606022
$values[5,13] = 555.924295494542
$values[5,14] = 0.005788675528669854
$values[5,15] = 0.005788675528669854

# Row 8: E,F,G,H,I,J,K,L,M,N,O,P,Q,R,S,T
$values[6,0] = 3
$values[6,1] = 1
$values[6,2] = 5.848171333333333
$values[6,3] = 17.544514
$values[6,4] = 0.02178075750932496
$values[6,5] = 0.02178075750932496
$values[6,6] = 3
$values[6,7] = 1
$values[6,8] = 14.90560833333333
$values[6,9] = 44.716825
$values[6,10] = 0.3750618357072762
$values[6,11] = 0.3750618357072762
$values[6,12] = 87.17055136089444
$values[6,13] = 784.53496224805
$values[6,14] = 0.008169130894542461
$values[6,15] = 0.008169130894542461

# Row 9: E,F,G,H,I,J,K,L,M,N,O,P,Q,R,S,T
$values[7,0] = 3
$values[7,1] = 1
$values[7,2] = 5.848171333333333
$values[7,3] = 17.544514
$values[7,4] = 0.02178075750932496
$values[7,5] = 0.02178075750932496
$values[7,6] = 3
$values[7,7] = 1
$values[7,8] = 3.885304333333333
$values[7,9] = 11.655913
$values[7,10] = 0.09776383109991163
$values[7,11] = 0.09776383109991162
$values[7,12] = 22.72192542347578
$values[7,13] = 204.497328811282
$values[7,14] = 0.002129370298369777
$values[7,15] = 0.002129370298369777

# Row 10: E,F,G,H,I,J,K,L,M,N,O,P,Q,R,S,T
$values[8,0] = 3
$values[8,1] = 1
$values[8,2] = 240.2954863333333
$values[8,3] = 720.8864590000001
$values[8,4] = 0.894949450023804
$values[8,5] = 0.894949450023804
$values[8,6] = 3
$values[8,7] = 1
$values[8,8] = 10.38865533333333
$values[8,9] = 31.165966
$values[8,10] = 0.261404167660619
$values[8,11] = 0.261404167660619
$values[8,12] = 2496.346985672711
$values[8,13] = 22467.1228710544
$values[8,14] = 0.2339435160818012
$values[8,15] = 0.2339435160818012

# Row 11: E,F,G,H,I,J,K,L,M,N,O,P,Q,R,S,T
$values[9,0] = 3
$values[9,1] = 1
$values[9,2] = 240.2954863333333
$values[9,3] = 720.8864590000001
$values[9,4] = 0.894949450023804
$values[9,5] = 0.894949450023804
$values[9,6] = 3
$values[9,7] = 1
$values[9,8] = 10.56216766666667
$values[9,9] = 31.686503
$values[9,10] = 0.2657701655321932
$values[9,11] = 0.2657701655321932
$values[9,12] = 2538.041216195875
$values[9,13] = 22842.37094576288
$values[9,14] = 0.2378508634757716
$values[9,15] = 0.2378508634757716

# Row 12: E,F,G,H,I,J,K,L,M,N,O,P,Q,R,S,T
$values[10,0] = 3
$values[10,1] = 1
$values[10,2] = 240.2954863333333
$values[10,3] = 720.8864590000001
$values[10,4] = 0.894949450023804
$values[10,5] = 0.894949450023804
$values[10,6] = 3
$values[10,7] = 1
$values[10,8] = 14.90560833333333
$values[10,9] = 44.716825
$values[10,10] = 0.3750618357072762
$values[10,11] = 0.3750618357072762
$values[10,12] = 3581.750403552519
$values[10,13] = 32235.75363197268
$values[10,14] = 0.3356613835911452
$values[10,15] = 0.3356613835911452

# Row 13: E,F,G,H,I,J,K,L,M,N,O,P,Q,R,S,T
$values[11,0] = 3
$values[11,1] = 1
$values[11,2] = 240.2954863333333
$values[11,3] = 720.8864590000001
$values[11,4] = 0.894949450023804
$values[11,5] = 0.894949450023804
$values[11,6] = 3
$values[11,7] = 1
$values[11,8] = 3.885304333333333
$values[11,9] = 11.655913
$values[11,10] = 0.09776383109991163
$values[11,11] = 0.09776383109991162
$values[11,12] = 933.6210943313408
$values[11,13] = 8402.589848982068
$values[11,14] = 0.08749368687508598
$values[11,15] = 0.08749368687508596

# Row 14: E,F,G,H,I,J,K,L,M,N,O,P,Q,R,S,T
$values[12,0] = 3
$values[12,1] = 1
$values[12,2] = 19.17634566666667
$values[12,3] = 57.529037
$values[12,4] = 0.07141981844820457
$values[12,5] = 0.07141981844820458
$values[12,6] = 3
$values[12,7] = 1
$values[12,8] = 10.38865533333333
$values[12,9] = 31.165966
$values[12,10] = 0.261404167660619
$values[12,11] = 0.261404167660619
$values[12,12] = 199.2164456838602
$values[12,13] = 1792.948011154742
$values[12,14] = 0.01866943819592544
$values[12,15] = 0.01866943819592544

# Row 15: E,F,G,H,I,J,K,L,M,N,O,P,Q,R,S,T
$values[13,0] = 3
$values[13,1] = 1
$values[13,2] = 19.17634566666667
$values[13,3] = 57.529037
$values[13,4] = 0.07141981844820457
$values[13,5] = 0.07141981844820458
$values[13,6] = 3
$values[13,7] = 1
$values[13,8] = 10.56216766666667
$values[13,9] = 31.686503
$values[13,10] = 0.2657701655321932
$values[13,11] = 0.2657701655321932
$values[13,12] = 202.5437781652901
$values[13,13] = 1822.894003487611
$values[13,14] = 0.01898125697125851
$values[13,15] = 0.01898125697125852

# Row 16: E,F,G,H,I,J,K,L,M,N,O,P,Q,R,S,T
$values[14,0] = 3
$values[14,1] = 1
$values[14,2] = 19.17634566666667
$values[14,3] = 57.529037
$values[14,4] = 0.07141981844820457
$values[14,5] = 0.07141981844820458
$values[14,6] = 3
$values[14,7] = 1
$values[14,8] = 14.90560833333333
$values[14,9] = 44.716825
$values[14,10] = 0.3750618357072762
$values[14,11] = 0.3750618357072762
$values[14,12] = 285.8350977719472
$values[14,13] = 2572.515879947525
$values[14,14] = 0.026786848213064
$values[14,15] = 0.026786848213064

# Row 17: E,F,G,H,I,J,K,L,M,N,O,P,Q,R,S,T
$values[15,0] = 3
$values[15,1] = 1
$values[15,2] = 19.17634566666667
$values[15,3] = 57.529037
$values[15,4] = 0.07141981844820457
$values[15,5] = 0.07141981844820458
$values[15,6] = 3
$values[15,7] = 1
$values[15,8] = 3.885304333333333
$values[15,9] = 11.655913
$values[15,10] = 0.09776383109991163
$values[15,11] = 0.09776383109991163
$values[15,12] = 74.50593891619789
$values[15,13] = 670.5534502457811
$values[15,14] = 0.006982275067956624
$values[15,15] = 0.006982275067956624

$ws.Range("E2:T17").Value = $values
